$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute(
        $findText,   # FindText
        $true,       # MatchCase
        $false,      # MatchWholeWord
        $false,      # MatchWildcards
        $false,      # MatchSoundsLike
        $false,      # MatchAllWordForms
        $true,       # Forward
        1,           # Wrap (wdFindContinue)
        $false,      # Format
        $replaceText,# ReplaceWith
        2            # Replace (wdReplaceAll)
    )
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) ", " + "soul" + ")"  ->  ", soul)"
Replace-Text ", soul)" ", soul)"

# 2) ", " + "'" + "of a" + " woman" + " in travail" + "'" + ") "  ->  ", 'of a woman in travail') "
Replace-Text ", ‘of a woman in travail’) " ", ‘of a woman in travail’) "

# 3) "'" + "spirit'"  ->  "'spirit'"   (disambiguated via surrounding context)
Replace-Text ", ‘spirit’) and" ", ‘spirit’) and"

# 4) "'" + "soul'"  ->  "'soul'"   (disambiguated via surrounding context)
Replace-Text ", ‘soul’) come from" ", ‘soul’) come from"

# 5) "evil spirit" + " or "  ->  "evil spirit or "
Replace-Text "evil spirit or " "evil spirit or "

# 6) "a " + "pleasing" + " " + "woman" + " later on"  ->  "a woman" + " " + "companion" + " later on"
Replace-Text "a pleasing woman later on" "a woman companion later on"
